# Fix a typo in the D100 magic item table on Feuil1: the range label for
# "Relancer sur table 8" was "96-00" (should be "96-100"), and move the
# selection cursor to where the author was last working (F21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Correct the "96-00" -> "96-100" typo (trailing space preserved, as in the
# original cell).
$ws.Range("A23").Value = "96-100 "

# Update the active selection to match the author's last cursor position.
$ws.Activate()
$ws.Range("F21").Select()
